$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Occasionally employed.deja.pro"
$ws.Range("C1").Value = "Regularly employed.deja.pro"
$ws.Range("D1").Value = "Student.deja.pro"
$ws.Range("E1").Value = "Unemployed / discouraged.deja.pro"
$ws.Range("F1").Value = "Receiving social benefits / pensioners / house-makers / disable.deja.pro"
$ws.Range("G1").Value = "Other.deja.pro"
$ws.Range("H1").Value = "Not known / missing.deja.pro"
$ws.Range("I1").Value = "Total.deja.pro"
